$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (price, volume %, hour) stay as text,
# matching the source data which stores every data cell as a string.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '256.49'
$ws.Range("E2").Value = '-0.33%'
$ws.Range("G2").Value = '17'

# Row 3
$ws.Range("D3").Value = '26.95'
$ws.Range("E3").Value = '-3.79%'
$ws.Range("G3").Value = '17'

# Row 4
$ws.Range("D4").Value = '4.642'
$ws.Range("E4").Value = '-10.98%'
$ws.Range("G4").Value = '17'

# Row 5
$ws.Range("D5").Value = '0.05899'
$ws.Range("E5").Value = '-0.09%'
$ws.Range("G5").Value = '17'

# Row 6
$ws.Range("D6").Value = '6.633'
$ws.Range("E6").Value = '-1.12%'
$ws.Range("G6").Value = '17'

# Row 7
$ws.Range("D7").Value = '0.8618'
$ws.Range("E7").Value = '-1.49%'
$ws.Range("G7").Value = '17'

# Row 8
$ws.Range("D8").Value = '0.9322'
$ws.Range("E8").Value = '-6.73%'
$ws.Range("G8").Value = '17'

# Row 9
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '0.1406'
$ws.Range("E9").Value = '-0.35%'
$ws.Range("G9").Value = '17'

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.03730'
$ws.Range("E10").Value = '2.36%'
$ws.Range("G10").Value = '17'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.07093'
$ws.Range("E11").Value = '-1.54%'
$ws.Range("G11").Value = '17'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.03230'
$ws.Range("E12").Value = '2.26%'
$ws.Range("G12").Value = '17'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09207'
$ws.Range("E13").Value = '-0.12%'
$ws.Range("G13").Value = '17'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001536'
$ws.Range("E14").Value = '-0.96%'
$ws.Range("G14").Value = '17'

# Row 15
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = '0.0006070'
$ws.Range("E15").Value = '-0.13%'
$ws.Range("G15").Value = '17'

# Row 16
$ws.Range("D16").Value = '0.006071'
$ws.Range("E16").Value = '3.69%'
$ws.Range("G16").Value = '17'

# Row 17
$ws.Range("D17").Value = '3.515'
$ws.Range("E17").Value = '0.62%'
$ws.Range("G17").Value = '17'

# Row 18
$ws.Range("D18").Value = '3.196'
$ws.Range("E18").Value = '-1.02%'
$ws.Range("G18").Value = '17'

# Row 19
$ws.Range("D19").Value = '2.223'
$ws.Range("E19").Value = '-0.12%'
$ws.Range("G19").Value = '17'

# Row 20
$ws.Range("E20").Value = '-0.74%'
$ws.Range("G20").Value = '17'

# Row 21
$ws.Range("D21").Value = '0.1267'
$ws.Range("E21").Value = '-1.75%'
$ws.Range("G21").Value = '17'

# Row 22
$ws.Range("D22").Value = '3.853'
$ws.Range("E22").Value = '9.31%'
$ws.Range("G22").Value = '17'

# Row 23
$ws.Range("D23").Value = '0.04214'
$ws.Range("E23").Value = '0.32%'
$ws.Range("G23").Value = '17'

# Row 24
$ws.Range("D24").Value = '0.001225'
$ws.Range("E24").Value = '0.62%'
$ws.Range("G24").Value = '17'

# Row 25
$ws.Range("D25").Value = '0.004279'
$ws.Range("E25").Value = '-5.92%'
$ws.Range("G25").Value = '17'

# Row 26
$ws.Range("E26").Value = '0.02%'
$ws.Range("G26").Value = '17'

# Row 27
$ws.Range("E27").Value = '-0.03%'
$ws.Range("G27").Value = '17'

# Row 28
$ws.Range("G28").Value = '17'

# Row 29
$ws.Range("G29").Value = '17'

# Row 30
$ws.Range("G30").Value = '17'

# Row 31
$ws.Range("G31").Value = '17'

# Row 32
$ws.Range("G32").Value = '17'

# Row 33
$ws.Range("G33").Value = '17'

# Row 34
$ws.Range("G34").Value = '17'

# Row 35
$ws.Range("G35").Value = '17'

# Row 36
$ws.Range("G36").Value = '17'

# Row 37
$ws.Range("G37").Value = '17'

# Row 38
$ws.Range("G38").Value = '17'

# Row 39
$ws.Range("G39").Value = '17'

# Row 40
$ws.Range("D40").Value = '0.03829'
$ws.Range("E40").Value = '-0.19%'
$ws.Range("G40").Value = '17'

# Row 41
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '0.1100'
$ws.Range("E41").Value = '-0.40%'
$ws.Range("G41").Value = '17'

# Row 42
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").Value = '0.003949'
$ws.Range("E42").Value = '-27.25%'
$ws.Range("G42").Value = '17'

# Row 43
$ws.Range("D43").Value = '0.002410'
$ws.Range("E43").Value = '-1.61%'
$ws.Range("G43").Value = '17'

# Row 44
$ws.Range("D44").Value = '0.01141'
$ws.Range("E44").Value = '7.06%'
$ws.Range("G44").Value = '17'

# Row 45
$ws.Range("D45").Value = '0.00005438'
$ws.Range("E45").Value = '0.31%'
$ws.Range("G45").Value = '17'

# Row 46
$ws.Range("E46").Value = '0.04%'
$ws.Range("G46").Value = '17'

# Row 47
$ws.Range("D47").Value = '0.06020'
$ws.Range("E47").Value = '-29.57%'
$ws.Range("G47").Value = '17'

# Row 48
$ws.Range("D48").Value = '0.002278'
$ws.Range("E48").Value = '6.60%'
$ws.Range("G48").Value = '17'

# Row 49
$ws.Range("E49").Value = '0.04%'
$ws.Range("G49").Value = '17'

# Row 50
$ws.Range("E50").Value = '0.04%'
$ws.Range("G50").Value = '17'

# Row 51
$ws.Range("G51").Value = '17'

